# Scheduled data refresh: update Leve price/profit columns (H-N)
# across the ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets with freshly pulled market data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4169288.5
$ws.Range("I32").Value = 1001.45
$ws.Range("K32").Value = 1001.45
$ws.Range("M32").Value = -714.45

$ws.Range("H74").Value = 4085.125
$ws.Range("I74").Value = 4041.7
$ws.Range("J74").Value = 4302.25
$ws.Range("K74").Value = 4041.7
$ws.Range("L74").Value = 4302.25
$ws.Range("M74").Value = -3167.7
$ws.Range("N74").Value = -6050.25

$ws.Range("H77").Value = 4085.125
$ws.Range("I77").Value = 4041.7
$ws.Range("J77").Value = 4302.25
$ws.Range("K77").Value = 20208.5
$ws.Range("L77").Value = 21511.25
$ws.Range("M77").Value = -15840.5
$ws.Range("N77").Value = -30247.25

$ws.Range("H132").Value = 1780.7046
$ws.Range("I132").Value = 1796.3658
$ws.Range("K132").Value = 5389.097400000001
$ws.Range("M132").Value = -2859.097400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3792.2144
$ws.Range("I86").Value = 1465.4445
$ws.Range("K86").Value = 1465.4445
$ws.Range("M86").Value = -342.4445000000001

$ws.Range("H89").Value = 3792.2144
$ws.Range("I89").Value = 1465.4445
$ws.Range("K89").Value = 7327.2225
$ws.Range("M89").Value = -1711.2225

$ws.Range("H134").Value = 1382.08
$ws.Range("I134").Value = 1202.2084
$ws.Range("K134").Value = 3606.6252
$ws.Range("M134").Value = -1071.6252

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3850.25
$ws.Range("I31").Value = 2440.6
$ws.Range("J31").Value = 4857.143
$ws.Range("K31").Value = 2440.6
$ws.Range("L31").Value = 4857.143
$ws.Range("M31").Value = -2145.6
$ws.Range("N31").Value = -5447.143

$ws.Range("H34").Value = 3850.25
$ws.Range("I34").Value = 2440.6
$ws.Range("J34").Value = 4857.143
$ws.Range("K34").Value = 2440.6
$ws.Range("L34").Value = 4857.143
$ws.Range("M34").Value = -2238.6
$ws.Range("N34").Value = -5261.143

$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").ClearContents()

$ws.Range("H94").Value = 3973.9
$ws.Range("I94").Value = 1176.5
$ws.Range("J94").Value = 5838.8335
$ws.Range("K94").Value = 1176.5
$ws.Range("L94").Value = 5838.8335
$ws.Range("M94").Value = -725.5
$ws.Range("N94").Value = -6740.8335

$ws.Range("H105").Value = 1319.4546
$ws.Range("I105").Value = 1167
$ws.Range("K105").Value = 1167
$ws.Range("M105").Value = 580

$ws.Range("H107").Value = 879.1818
$ws.Range("I107").Value = 774.06665
$ws.Range("J107").Value = 1104.4286
$ws.Range("K107").Value = 774.06665
$ws.Range("L107").Value = 1104.4286
$ws.Range("M107").Value = 1145.93335
$ws.Range("N107").Value = -4944.4286

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2262.25
$ws.Range("I5").Value = 2066
$ws.Range("J5").Value = 2380
$ws.Range("K5").Value = 6198
$ws.Range("L5").Value = 7140
$ws.Range("M5").Value = -6086
$ws.Range("N5").Value = -7364

$ws.Range("H47").Value = 863.25
$ws.Range("J47").Value = 725
$ws.Range("L47").Value = 2175
$ws.Range("N47").Value = -3037

$ws.Range("H55").Value = 2701.5
$ws.Range("J55").Value = 4374
$ws.Range("L55").Value = 13122
$ws.Range("N55").Value = -13476

$ws.Range("H135").Value = 2262.25
$ws.Range("I135").Value = 2066
$ws.Range("J135").Value = 2380
$ws.Range("K135").Value = 18594
$ws.Range("L135").Value = 21420
$ws.Range("M135").Value = -16059
$ws.Range("N135").Value = -26490

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 39591
$ws.Range("I69").Value = 39591
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 39591
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -38842
$ws.Range("N69").ClearContents()

$ws.Range("H72").Value = 39591
$ws.Range("I72").Value = 39591
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 118773
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -115029
$ws.Range("N72").ClearContents()

$ws.Range("H132").Value = 38925.93
$ws.Range("I132").Value = 49756.145
$ws.Range("J132").Value = 6435.2856
$ws.Range("K132").Value = 149268.435
$ws.Range("L132").Value = 19305.8568
$ws.Range("M132").Value = -146738.435
$ws.Range("N132").Value = -24365.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()

$ws.Range("H55").Value = 1188.5
$ws.Range("I55").Value = 1343.9
$ws.Range("J55").Value = 800
$ws.Range("K55").Value = 1343.9
$ws.Range("L55").Value = 800
$ws.Range("M55").Value = -1170.9
$ws.Range("N55").Value = -1146

$ws.Range("H132").Value = 10650
$ws.Range("I132").Value = 7223.75
$ws.Range("J132").Value = 17502.5
$ws.Range("K132").Value = 21671.25
$ws.Range("L132").Value = 52507.5
$ws.Range("M132").Value = -19141.25
$ws.Range("N132").Value = -57567.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 5000000
$ws.Range("I21").Value = 5000000
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 5000000
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -4999765
$ws.Range("N21").ClearContents()

$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()

$ws.Range("H28").Value = 21000
$ws.Range("J28").Value = 21000
$ws.Range("L28").Value = 21000
$ws.Range("N28").Value = -21696

$ws.Range("H35").Value = 5000000
$ws.Range("I35").Value = 5000000
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 5000000
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -4999710
$ws.Range("N35").ClearContents()

$ws.Range("H126").Value = 2814.4285
$ws.Range("I126").Value = 1279.1578
$ws.Range("J126").Value = 6055.5557
$ws.Range("K126").Value = 3837.4734
$ws.Range("L126").Value = 18166.6671
$ws.Range("M126").Value = -1367.4734
$ws.Range("N126").Value = -23106.6671

$ws.Range("H136").Value = 3863.4707
$ws.Range("I136").Value = 3512.8572
$ws.Range("J136").Value = 5499.6665
$ws.Range("K136").Value = 10538.5716
$ws.Range("L136").Value = 16498.9995
$ws.Range("M136").Value = -7988.571599999999
$ws.Range("N136").Value = -21598.9995

